$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the F:V payload between row 95 and row 96 (columns A:E - Indice,
#    pais, torneio, temporada, data_partida - stay exactly as they are).
# ---------------------------------------------------------------------------

# NOTE: `.Value` reads back a wrapper/descriptor instead of the scalar in
# this host, so all reads below go through `.Value2` (which yields the
# real scalar for both text and numeric cells, without reinterpreting the
# "dd/mm/yyyy hh:mm" text columns as dates).

$row95 = @{
    F = $ws.Cells.Item(95, 6).Value2
    G = $ws.Cells.Item(95, 7).Value2
    H = $ws.Cells.Item(95, 8).Value2
    I = $ws.Cells.Item(95, 9).Value2
    J = $ws.Cells.Item(95, 10).Value2
    K = $ws.Cells.Item(95, 11).Value2
    L = $ws.Cells.Item(95, 12).Value2
    M = $ws.Cells.Item(95, 13).Value2
    N = $ws.Cells.Item(95, 14).Value2
    O = $ws.Cells.Item(95, 15).Value2
    P = $ws.Cells.Item(95, 16).Value2
    Q = $ws.Cells.Item(95, 17).Value2
    R = $ws.Cells.Item(95, 18).Value2
    S = $ws.Cells.Item(95, 19).Value2
    T = $ws.Cells.Item(95, 20).Value2
    U = $ws.Cells.Item(95, 21).Value2
    V = $ws.Cells.Item(95, 22).Value2
}

$row96 = @{
    F = $ws.Cells.Item(96, 6).Value2
    G = $ws.Cells.Item(96, 7).Value2
    H = $ws.Cells.Item(96, 8).Value2
    I = $ws.Cells.Item(96, 9).Value2
    J = $ws.Cells.Item(96, 10).Value2
    K = $ws.Cells.Item(96, 11).Value2
    L = $ws.Cells.Item(96, 12).Value2
    M = $ws.Cells.Item(96, 13).Value2
    N = $ws.Cells.Item(96, 14).Value2
    O = $ws.Cells.Item(96, 15).Value2
    P = $ws.Cells.Item(96, 16).Value2
    Q = $ws.Cells.Item(96, 17).Value2
    R = $ws.Cells.Item(96, 18).Value2
    S = $ws.Cells.Item(96, 19).Value2
    T = $ws.Cells.Item(96, 20).Value2
    U = $ws.Cells.Item(96, 21).Value2
    V = $ws.Cells.Item(96, 22).Value2
}

$ws.Cells.Item(95, 6).Value2 = $row96.F
$ws.Cells.Item(95, 7).Value2 = $row96.G
$ws.Cells.Item(95, 8).Value2 = $row96.H
$ws.Cells.Item(95, 9).Value2 = $row96.I
$ws.Cells.Item(95, 10).Value2 = $row96.J
$ws.Cells.Item(95, 11).Value2 = $row96.K
$ws.Cells.Item(95, 12).Value2 = $row96.L
$ws.Cells.Item(95, 13).Value2 = $row96.M
$ws.Cells.Item(95, 14).Value2 = $row96.N
$ws.Cells.Item(95, 15).Value2 = $row96.O
$ws.Cells.Item(95, 16).Value2 = $row96.P
$ws.Cells.Item(95, 17).Value2 = $row96.Q
$ws.Cells.Item(95, 18).Value2 = $row96.R
$ws.Cells.Item(95, 19).Value2 = $row96.S
$ws.Cells.Item(95, 20).Value2 = $row96.T
$ws.Cells.Item(95, 21).Value2 = $row96.U
$ws.Cells.Item(95, 22).Value2 = $row96.V

$ws.Cells.Item(96, 6).Value2 = $row95.F
$ws.Cells.Item(96, 7).Value2 = $row95.G
$ws.Cells.Item(96, 8).Value2 = $row95.H
$ws.Cells.Item(96, 9).Value2 = $row95.I
$ws.Cells.Item(96, 10).Value2 = $row95.J
$ws.Cells.Item(96, 11).Value2 = $row95.K
$ws.Cells.Item(96, 12).Value2 = $row95.L
$ws.Cells.Item(96, 13).Value2 = $row95.M
$ws.Cells.Item(96, 14).Value2 = $row95.N
$ws.Cells.Item(96, 15).Value2 = $row95.O
$ws.Cells.Item(96, 16).Value2 = $row95.P
$ws.Cells.Item(96, 17).Value2 = $row95.Q
$ws.Cells.Item(96, 18).Value2 = $row95.R
$ws.Cells.Item(96, 19).Value2 = $row95.S
$ws.Cells.Item(96, 20).Value2 = $row95.T
$ws.Cells.Item(96, 21).Value2 = $row95.U
$ws.Cells.Item(96, 22).Value2 = $row95.V

# ---------------------------------------------------------------------------
# 2) Append two new match rows (119, 120) at the bottom of the sheet,
#    copying the formatting of the last existing row (118) first so the
#    bold/bordered index column and the custom date-time column keep
#    their original styles (s="1" / s="2").
# ---------------------------------------------------------------------------

$ws.Range("A118:V118").Copy()
$ws.Range("A119:V120").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 119
$ws.Cells.Item(119, 1).Value = 118
$ws.Cells.Item(119, 2).Value = "costa-rica"
$ws.Cells.Item(119, 3).Value = "primera-division"
$ws.Cells.Item(119, 4).Value = "2023-2024"
$ws.Cells.Item(119, 5).Value = 45242.95833333334
$ws.Cells.Item(119, 6).Value = "Saprissa"
$ws.Cells.Item(119, 7).Value = 4
$ws.Cells.Item(119, 8).Value = "Zeledon"
$ws.Cells.Item(119, 9).Value = 0
$ws.Cells.Item(119, 10).Value = 1.15
$ws.Cells.Item(119, 11).Value = "10/11/2023 03:13"
$ws.Cells.Item(119, 12).Value = 1.18
$ws.Cells.Item(119, 13).Value = "12/11/2023 22:59"
$ws.Cells.Item(119, 14).Value = 7.29
$ws.Cells.Item(119, 15).Value = "10/11/2023 03:13"
$ws.Cells.Item(119, 16).Value = 7.52
$ws.Cells.Item(119, 17).Value = "12/11/2023 22:59"
$ws.Cells.Item(119, 18).Value = 12.6
$ws.Cells.Item(119, 19).Value = "10/11/2023 03:13"
$ws.Cells.Item(119, 20).Value = 14.05
$ws.Cells.Item(119, 21).Value = "12/11/2023 22:59"
$ws.Cells.Item(119, 22).Value = "https://www.betexplorer.com/football/costa-rica/primera-division/saprissa-zeledon/GhIsFOnc/"

# Row 120
$ws.Cells.Item(120, 1).Value = 119
$ws.Cells.Item(120, 2).Value = "costa-rica"
$ws.Cells.Item(120, 3).Value = "primera-division"
$ws.Cells.Item(120, 4).Value = "2023-2024"
$ws.Cells.Item(120, 5).Value = 45243.04166666666
$ws.Cells.Item(120, 6).Value = "Sporting San Jose"
$ws.Cells.Item(120, 7).Value = 1
$ws.Cells.Item(120, 8).Value = "AD Santos"
$ws.Cells.Item(120, 9).Value = 0
$ws.Cells.Item(120, 10).Value = 2
$ws.Cells.Item(120, 11).Value = "09/11/2023 03:42"
$ws.Cells.Item(120, 12).Value = 2.12
$ws.Cells.Item(120, 13).Value = "13/11/2023 00:56"
$ws.Cells.Item(120, 14).Value = 3.43
$ws.Cells.Item(120, 15).Value = "09/11/2023 03:42"
$ws.Cells.Item(120, 16).Value = 3.4
$ws.Cells.Item(120, 17).Value = "13/11/2023 00:56"
$ws.Cells.Item(120, 18).Value = 3.55
$ws.Cells.Item(120, 19).Value = "09/11/2023 03:42"
$ws.Cells.Item(120, 20).Value = 3.58
$ws.Cells.Item(120, 21).Value = "13/11/2023 00:56"
$ws.Cells.Item(120, 22).Value = "https://www.betexplorer.com/football/costa-rica/primera-division/sporting-san-jose-santos-de-guapiles/29GkD2H9/"

# ---------------------------------------------------------------------------
# 3) Keep the sheet dimension in sync with the new extent.
# ---------------------------------------------------------------------------
$ws.UsedRange | Out-Null
